$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign people to teams: Team 1 grows to 6 members, Team 2 shrinks to 5.
# New order (rows 2-12): Veselin, Rawda, Hannah, Mirit, Bogdana, Martin (Team 1)
#                         Alice, Bob, Claire, David, Elaine (Team 2)
$names = @("Veselin", "Rawda", "Hannah", "Mirit", "Bogdana", "Martin", "Alice", "Bob", "Claire", "David", "Elaine")
$teams = @(1, 1, 1, 1, 1, 1, 1, 2, 2, 2, 2)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $teams[$i]
}

# Update the selected cell/range to match the saved view state.
$ws.Range("C10").Select()
